$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update G1 header text (camp -> camp with legend)
$ws.Range("G1").Value = "阵营1=默认，2=玩家，3=怪"

# Header row style tweaks: G1 gains wrap-text style; J1:S1 keep wrap-text (style index changes in
# the source file, but both indices render identically - wrapText + vertical-center)
$ws.Range("G1").WrapText = $true
$ws.Range("J1:S1").WrapText = $true

# Column G width change (10.375 -> ~12.125; engine quantizes to 1/7 steps)
$ws.Columns.Item(7).ColumnWidth = 11.4

# Populate / overwrite NPC data rows 4-23
# Row 4
$ws.Range("A4").Value = 100010001
$ws.Range("B4").Value = "刀盾兵"
$ws.Range("C4").Value = 1000016
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("L4").Value = 2
$ws.Range("M4").Value = "(Att,2),(Hp,10)"
$ws.Range("N4").Value = "(Speed,5),(Vision,10)"
$ws.Range("Q4").Value = 2
# Row 5
$ws.Range("A5").Value = 100020001
$ws.Range("B5").Value = "巨盾兵"
$ws.Range("C5").Value = 1000017
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = 2
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("L5").Value = 3
$ws.Range("M5").Value = "(Att,2),(Hp,20)"
$ws.Range("N5").Value = "(Speed,5),(Vision,10)"
$ws.Range("Q5").Value = 5
# Row 6
$ws.Range("A6").Value = 100030001
$ws.Range("B6").Value = "短弓手"
$ws.Range("C6").Value = 1000018
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 2
$ws.Range("L6").Value = 2
$ws.Range("M6").Value = "(Att,4),(Hp,5)"
$ws.Range("N6").Value = "(Speed,5),(Vision,10)"
$ws.Range("Q6").Value = 2
# Row 7
$ws.Range("A7").Value = 100040001
$ws.Range("B7").Value = "长弓手"
$ws.Range("C7").Value = 1000019
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 4
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 2
$ws.Range("L7").Value = 3
$ws.Range("M7").Value = "(Att,5),(Hp,8)"
$ws.Range("N7").Value = "(Speed,5),(Vision,10)"
$ws.Range("Q7").Value = 5
# Row 8
$ws.Range("A8").Value = 100050001
$ws.Range("B8").Value = "轻骑兵"
$ws.Range("C8").Value = 1000020
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 3
$ws.Range("L8").Value = 2
$ws.Range("M8").Value = "(Att,3),(Hp,7)"
$ws.Range("N8").Value = "(Speed,10),(Vision,10)"
$ws.Range("Q8").Value = 2
# Row 9
$ws.Range("A9").Value = 100060001
$ws.Range("B9").Value = "重骑兵"
$ws.Range("C9").Value = 1000021
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = 2
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 3
$ws.Range("L9").Value = 3
$ws.Range("M9").Value = "(Att,4),(Hp,10)"
$ws.Range("N9").Value = "(Speed,10),(Vision,10)"
$ws.Range("Q9").Value = 5
# Row 10
$ws.Range("A10").Value = 100070001
$ws.Range("B10").Value = "冲车"
$ws.Range("C10").Value = 1000022
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 3
$ws.Range("G10").Value = 2
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 4
$ws.Range("L10").Value = 2
$ws.Range("M10").Value = "(Att,1),(Hp,20)"
$ws.Range("N10").Value = "(Speed,3),(Vision,10)"
$ws.Range("Q10").Value = 2
# Row 11
$ws.Range("A11").Value = 100080001
$ws.Range("B11").Value = "投石车"
$ws.Range("C11").Value = 1000023
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 4
$ws.Range("G11").Value = 2
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 4
$ws.Range("L11").Value = 3
$ws.Range("M11").Value = "(Att,2),(Hp,30)"
$ws.Range("N11").Value = "(Speed,3),(Vision,10)"
$ws.Range("Q11").Value = 5
# Row 12
$ws.Range("A12").Value = 100090001
$ws.Range("B12").Value = "军医"
$ws.Range("C12").Value = 1000024
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2
$ws.Range("H12").Value = 1
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 5
$ws.Range("L12").Value = 2
$ws.Range("M12").Value = "(Att,1),(Hp,15)"
$ws.Range("N12").Value = "(Speed,5),(Vision,10)"
$ws.Range("Q12").Value = 2
# Row 13
$ws.Range("A13").Value = 100100001
$ws.Range("B13").Value = "巫医"
$ws.Range("C13").Value = 1000025
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 5
$ws.Range("L13").Value = 3
$ws.Range("M13").Value = "(Att,2),(Hp,25)"
$ws.Range("N13").Value = "(Speed,5),(Vision,10)"
$ws.Range("Q13").Value = 5
# Row 14
$ws.Range("A14").Value = 200010001
$ws.Range("B14").Value = "刀盾兵"
$ws.Range("C14").Value = 1000016
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 3
$ws.Range("H14").Value = 1
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 1
$ws.Range("L14").Value = 2
$ws.Range("M14").Value = "(Att,2),(Hp,10)"
$ws.Range("N14").Value = "(Speed,5),(Vision,10)"
$ws.Range("Q14").Value = 2
# Row 15
$ws.Range("A15").Value = 200020001
$ws.Range("B15").Value = "巨盾兵"
$ws.Range("C15").Value = 1000017
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 1
$ws.Range("L15").Value = 3
$ws.Range("M15").Value = "(Att,2),(Hp,20)"
$ws.Range("N15").Value = "(Speed,5),(Vision,10)"
$ws.Range("Q15").Value = 5
# Row 16
$ws.Range("A16").Value = 200030001
$ws.Range("B16").Value = "短弓手"
$ws.Range("C16").Value = 1000018
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 1
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 2
$ws.Range("L16").Value = 2
$ws.Range("M16").Value = "(Att,4),(Hp,5)"
$ws.Range("N16").Value = "(Speed,5),(Vision,10)"
$ws.Range("Q16").Value = 2
# Row 17
$ws.Range("A17").Value = 200040001
$ws.Range("B17").Value = "长弓手"
$ws.Range("C17").Value = 1000019
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 4
$ws.Range("G17").Value = 3
$ws.Range("H17").Value = 1
$ws.Range("I17").Value = 1
$ws.Range("J17").Value = 2
$ws.Range("L17").Value = 3
$ws.Range("M17").Value = "(Att,5),(Hp,8)"
$ws.Range("N17").Value = "(Speed,5),(Vision,10)"
$ws.Range("Q17").Value = 5
# Row 18
$ws.Range("A18").Value = 200050001
$ws.Range("B18").Value = "轻骑兵"
$ws.Range("C18").Value = 1000020
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 1
$ws.Range("I18").Value = 1
$ws.Range("J18").Value = 3
$ws.Range("L18").Value = 2
$ws.Range("M18").Value = "(Att,3),(Hp,7)"
$ws.Range("N18").Value = "(Speed,10),(Vision,10)"
$ws.Range("Q18").Value = 2
# Row 19
$ws.Range("A19").Value = 200060001
$ws.Range("B19").Value = "重骑兵"
$ws.Range("C19").Value = 1000021
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 2
$ws.Range("G19").Value = 3
$ws.Range("H19").Value = 1
$ws.Range("I19").Value = 1
$ws.Range("J19").Value = 3
$ws.Range("L19").Value = 3
$ws.Range("M19").Value = "(Att,4),(Hp,10)"
$ws.Range("N19").Value = "(Speed,10),(Vision,10)"
$ws.Range("Q19").Value = 5
# Row 20
$ws.Range("A20").Value = 200070001
$ws.Range("B20").Value = "冲车"
$ws.Range("C20").Value = 1000022
$ws.Range("E20").Value = 2
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 1
$ws.Range("I20").Value = 1
$ws.Range("J20").Value = 4
$ws.Range("L20").Value = 2
$ws.Range("M20").Value = "(Att,1),(Hp,20)"
$ws.Range("N20").Value = "(Speed,3),(Vision,10)"
$ws.Range("Q20").Value = 2
# Row 21
$ws.Range("A21").Value = 200080001
$ws.Range("B21").Value = "投石车"
$ws.Range("C21").Value = 1000023
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 4
$ws.Range("G21").Value = 3
$ws.Range("H21").Value = 1
$ws.Range("I21").Value = 1
$ws.Range("J21").Value = 4
$ws.Range("L21").Value = 3
$ws.Range("M21").Value = "(Att,2),(Hp,30)"
$ws.Range("N21").Value = "(Speed,3),(Vision,10)"
$ws.Range("Q21").Value = 5
# Row 22
$ws.Range("A22").Value = 200090001
$ws.Range("B22").Value = "军医"
$ws.Range("C22").Value = 1000024
$ws.Range("E22").Value = 2
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 1
$ws.Range("I22").Value = 1
$ws.Range("J22").Value = 5
$ws.Range("L22").Value = 2
$ws.Range("M22").Value = "(Att,1),(Hp,15)"
$ws.Range("N22").Value = "(Speed,5),(Vision,10)"
$ws.Range("Q22").Value = 2
# Row 23
$ws.Range("A23").Value = 200100001
$ws.Range("B23").Value = "巫医"
$ws.Range("C23").Value = 1000025
$ws.Range("E23").Value = 2
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 1
$ws.Range("I23").Value = 1
$ws.Range("J23").Value = 5
$ws.Range("L23").Value = 3
$ws.Range("M23").Value = "(Att,2),(Hp,25)"
$ws.Range("N23").Value = "(Speed,5),(Vision,10)"
$ws.Range("Q23").Value = 5

# Update selection to match target (cosmetic)
$ws.Range("G33").Select()